$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume%) hold text-formatted numeric-looking strings
# (e.g. "8.00", "64.788.68") in the source data. Force the whole D2:D51 range
# to Text format before writing so Excel keeps these as literal strings instead
# of auto-converting them to numbers (which would drop formatting/precision),
# then restore the default style so no stray per-cell formatting is introduced.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '64.788.68'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').Value = '3.172.19'
$ws.Range('E3').Value = '  +1.47%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '616.91'
$ws.Range('E5').Value = '  +3.81%  '
$ws.Range('D6').Value = '148.03'
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('D8').Value = '3.169.36'
$ws.Range('E8').Value = '  +1.26%  '
$ws.Range('D9').Value = '0.531'
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('D11').Value = '5.50'
$ws.Range('E11').Value = '  -2.45%  '
$ws.Range('D12').Value = '0.475'
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').Value = '35.93'
$ws.Range('E14').Value = '  -2.35%  '
$ws.Range('D15').Value = '3.690.86'
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('E16').Value = '  +3.23%  '
$ws.Range('D17').Value = '64.796.99'
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').Value = '3.168.21'
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('D20').Value = '481.06'
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('D21').Value = '14.79'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').Value = '0.726'
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('D23').Value = '8.00'
$ws.Range('E23').Value = '  +3.26%  '
$ws.Range('D24').Value = '13.82'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').Value = '84.63'
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  -1.97%  '
$ws.Range('D28').Value = '8.61'
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').Value = '6.96'
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('D30').Value = '0.117'
$ws.Range('E30').Value = '  -5.21%  '
$ws.Range('E31').Value = '  -6.87%  '
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('D34').Value = '26.60'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('E35').Value = '  +2.27%  '
$ws.Range('D36').Value = '0.0₃0780'
$ws.Range('E36').Value = '  +5.02%  '
$ws.Range('D37').Value = '6.03'
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('D38').Value = '3.23'
$ws.Range('E38').Value = '  +0.68%  '
$ws.Range('D39').Value = '53.11'
$ws.Range('E39').Value = '  -2.81%  '
$ws.Range('D40').Value = '461.66'
$ws.Range('E40').Value = '  +3.01%  '
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('E42').Value = '  -2.58%  '
$ws.Range('D43').Value = '8.44'
$ws.Range('E43').Value = '  -0.67%  '
$ws.Range('D44').Value = '2.851.29'
$ws.Range('E44').Value = '  -0.97%  '
$ws.Range('E45').Value = '  -2.63%  '
$ws.Range('E46').Value = '  -1.22%  '
$ws.Range('D47').Value = '2.46'
$ws.Range('E47').Value = '  +6.06%  '
$ws.Range('D48').Value = '26.71'
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('E50').Value = '  -0.92%  '
$ws.Range('D51').Value = '120.52'
$ws.Range('E51').Value = '  +1.54%  '

$priceRange.Style = "Normal"
